# Apply odds updates for 2025-12-18 Betfair Back/Lay workbook
# (values taken from the canonical OOXML diff, cell-by-cell)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("R2").Value = 1.21
# Row 3
$ws.Range("Q3").Value = 1.86
# Row 4
$ws.Range("G4").Value = 1.26
$ws.Range("H4").Value = 15
$ws.Range("I4").Value = 15.5
$ws.Range("J4").Value = 7.4
$ws.Range("K4").Value = 7.6
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.17
$ws.Range("S4").Value = 2.3
$ws.Range("T4").Value = 2.28
$ws.Range("U4").Value = 1.69
$ws.Range("X4").Value = 34
$ws.Range("Y4").Value = 60
$ws.Range("Z4").Value = 190
$ws.Range("AB4").Value = 13
$ws.Range("AC4").Value = 22
$ws.Range("AD4").Value = 75
$ws.Range("AF4").Value = 9.199999999999999
$ws.Range("AG4").Value = 17.5
$ws.Range("AH4").Value = 50
$ws.Range("AI4").Value = 240
$ws.Range("AJ4").Value = 9.199999999999999
$ws.Range("AK4").Value = 15
$ws.Range("AL4").Value = 65
$ws.Range("AM4").Value = 270
# Row 5
$ws.Range("G5").Value = 1.69
$ws.Range("K5").Value = 4.4
# Row 6
$ws.Range("F6").Value = 1.92
$ws.Range("H6").Value = 4.4
$ws.Range("J6").Value = 3.8
# Row 7
$ws.Range("F7").Value = 1.29
$ws.Range("H7").Value = 12
$ws.Range("K7").Value = 6.8
$ws.Range("P7").Value = 2.36
$ws.Range("Q7").Value = 1.63
$ws.Range("T7").Value = 2.2
$ws.Range("U7").Value = 1.75
$ws.Range("X7").Value = 1000
$ws.Range("Z7").Value = 150
$ws.Range("AC7").Value = 14.5
$ws.Range("AE7").Value = 270
$ws.Range("AI7").Value = 210
$ws.Range("AM7").Value = 210
# Row 8
$ws.Range("F8").Value = 1.96
$ws.Range("I8").Value = 4.6
$ws.Range("J8").Value = 3.7
$ws.Range("K8").Value = 3.9
$ws.Range("AB8").Value = 9
$ws.Range("AD8").Value = 1000
# Row 9
$ws.Range("J9").Value = 3.3
$ws.Range("P9").Value = 1.81
$ws.Range("X9").Value = 11
# Row 10
$ws.Range("F10").Value = 1.83
$ws.Range("G10").Value = 1.99
$ws.Range("H10").Value = 4.3
$ws.Range("I10").Value = 5.1
$ws.Range("P10").Value = 1.79
$ws.Range("AB10").Value = 10.5
# Row 11
$ws.Range("G11").Value = 3.6
$ws.Range("H11").Value = 2.74
$ws.Range("I11").Value = 2.92
$ws.Range("P11").Value = 1.55
$ws.Range("T11").Value = 1.86
$ws.Range("U11").Value = 1.72
# Row 12
$ws.Range("F12").Value = 1.25
$ws.Range("G12").Value = 1.26
$ws.Range("H12").Value = 14.5
$ws.Range("I12").Value = 20
$ws.Range("J12").Value = 6.8
$ws.Range("K12").Value = 7.6
$ws.Range("P12").Value = 2.56
$ws.Range("Q12").Value = 1.52
$ws.Range("X12").Value = 34
$ws.Range("Z12").Value = 190
$ws.Range("AF12").Value = 8.800000000000001
$ws.Range("AH12").Value = 48
$ws.Range("AI12").Value = 240
$ws.Range("AJ12").Value = 10
$ws.Range("AK12").Value = 17
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 4.1
# Row 13
$ws.Range("K13").Value = 5.3
$ws.Range("Q13").Value = 1.84
$ws.Range("AH13").Value = 1000
# Row 14
$ws.Range("P14").Value = 2.64
# Row 15
$ws.Range("F15").Value = 2.22
$ws.Range("H15").Value = 3.65
$ws.Range("I15").Value = 3.7
$ws.Range("J15").Value = 3.6
$ws.Range("K15").Value = 3.65
# Row 16
$ws.Range("G16").Value = 1.32
$ws.Range("H16").Value = 12
$ws.Range("I16").Value = 14.5
$ws.Range("J16").Value = 6.4
$ws.Range("U16").Value = 1.92
$ws.Range("Y16").Value = 1000
# Row 17
$ws.Range("F17").Value = 3.5
$ws.Range("G17").Value = 3.55
# Row 18
$ws.Range("F18").Value = 1.61
$ws.Range("G18").Value = 1.64
$ws.Range("I18").Value = 6.6
$ws.Range("J18").Value = 4.4
$ws.Range("P18").Value = 2.4
$ws.Range("X18").Value = 26
$ws.Range("AC18").Value = 11.5
$ws.Range("AH18").Value = 20
$ws.Range("AK18").Value = 15.5
$ws.Range("AL18").Value = 29
$ws.Range("AN18").Value = 6.8
# Row 19
$ws.Range("F19").Value = 1.23
$ws.Range("G19").Value = 1.26
$ws.Range("I19").Value = 17
$ws.Range("J19").Value = 7.2
$ws.Range("P19").Value = 3.3
$ws.Range("Q19").Value = 1.37
$ws.Range("U19").Value = 2.02
# Row 20
$ws.Range("G20").Value = 1.73
$ws.Range("H20").Value = 5.1
$ws.Range("K20").Value = 5
$ws.Range("M20").Value = 1.04
$ws.Range("X20").Value = 25
$ws.Range("AA20").Value = 150
$ws.Range("AB20").Value = 14
